# Update specific values in column E per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E12" = 17.386
    "E32" = 17.4
    "E36" = 16.652
    "E38" = 16.631
    "E46" = 16.815
    "E54" = 16.87
    "E55" = 16.494
    "E67" = 17.24
    "E69" = 17.386
    "E72" = 16.939
    "E91" = 17.577
    "E99" = 16.659
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
